# Replace embedded line-breaks in a handful of shared-string cell values
# with single spaces (the author's "join the wrapped lines back into one
# line" cleanup pass across the Adult/Pediatric VFC and influenza sheets).

$wb = $excel.ActiveWorkbook

$wsAdultVFC  = $wb.Worksheets.Item("Adult VFC Vaccine ")
$wsPedFlu    = $wb.Worksheets.Item("Pediatric influenza Influenza")
$wsAdultFlu  = $wb.Worksheets.Item("Adult influenza VFC")

# Adult VFC Vaccine  sheet
$wsAdultVFC.Range("B15").Value = "Tetanus  Diphtheria Toxoids Adsorbed for Adults No Preservative"
$wsAdultVFC.Range("H15").Value = "MassBioLogics (Akorn, Inc)"

# Pediatric influenza Influenza sheet
$wsPedFlu.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$wsPedFlu.Range("B8").Value = "FluMist No Preservative"

# Adult influenza VFC sheet
$wsAdultFlu.Range("B10").Value = "Afluria No Preservative"
$wsAdultFlu.Range("H10").Value = "CSL Biotherapies"
$wsAdultFlu.Range("H11").Value = "CSL Biotherapies"
